# V-King PRO 3D 400 BOM - "Cleaning up BOM part names"
#
# The author trimmed Fusion-360-exported version suffixes (e.g. " v21",
# " v3", " HEX v7") from Part Name values in column A, normalized a few
# inconsistent capitalisations (e.g. "M5 LOCK NUT" -> "M5 Lock Nut") to
# match column B, filled in a handful of previously-blank Description
# cells (column C), and renamed two parts/descriptions that were
# mis-labelled (rows 31 "Bed Heater" and 32 "Build Plate").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Cabinet parts: fill in previously blank Description cells ---
$ws.Range("C8").Value  = "Filament Holder Rod"
$ws.Range("C9").Value  = "LCD Rear Cover"
$ws.Range("C10").Value = "LCD Front Cover"
$ws.Range("C11").Value = "Filament Holder Nut"

# --- Controller components: drop Fusion 360 version suffixes from Part Name ---
$ws.Range("A19").Value = "12864 LCD"
$ws.Range("A20").Value = "Filament Sensor "
$ws.Range("A21").Value = "SSR Heat Sink Fan "
$ws.Range("A22").Value = "FAN 60-10 "
$ws.Range("A23").Value = "FAN 40-10mm "
$ws.Range("A25").Value = "Optical Sensor "

# --- Heatbed components ---
$ws.Range("A30").Value = "ssr-40 DA "
$ws.Range("A31").Value = "Silicone Heatpad"
$ws.Range("B31").Value = "220VAC 1000W Silicone Pad"
$ws.Range("C31").Value = "220AC 1000w 400*400"
$ws.Range("A32").Value = "Aluminum Build Plate"
$ws.Range("B32").Value = "400*400*6 Build Plate"

# --- Mechatronical components ---
$ws.Range("A35").Value = "E3d Hotend "
$ws.Range("A41").Value = "GT2 Motor Pulley 20T "
$ws.Range("A42").Value = "GT2 Pulley 20NT "
$ws.Range("A43").Value = "GT2 Pulley 20T "

# --- Fastners and Bearings ---
$ws.Range("A46").Value = "M5 8OD Washer "
$ws.Range("B46").Value = "M5 8OD Washer"
$ws.Range("A47").Value = "Bearing 5-16-5 "
$ws.Range("A48").Value = "M5 T-Nut Square "
$ws.Range("A49").Value = "M5 Tee Nut "
$ws.Range("A50").Value = "M5x35 "
$ws.Range("A51").Value = "M5x30 "
$ws.Range("A52").Value = "M5x30 Hex "
$ws.Range("A53").Value = "M5x20 Hex "
$ws.Range("A54").Value = "M5X25 "
$ws.Range("A55").Value = "M5x10 "
$ws.Range("A56").Value = "M5x10 Hex "
$ws.Range("A57").Value = "M4x20 "
$ws.Range("A58").Value = "M4x16 "
$ws.Range("A59").Value = "M4x12 "
$ws.Range("A60").Value = "M4x10 "
$ws.Range("A61").Value = "M3X25 "
$ws.Range("A62").Value = "M5x8 V"
$ws.Range("A63").Value = "M5 Lock Nut "
$ws.Range("B63").Value = "M5 Lock Nut"
$ws.Range("A64").Value = "M3x20 "
$ws.Range("A65").Value = "M5 Nut "
$ws.Range("A66").Value = "M3x14 "
$ws.Range("A67").Value = "M3x12 "
$ws.Range("A68").Value = "M3x10 "
$ws.Range("A69").Value = "M4 NUT "
$ws.Range("A70").Value = "M3x8 "
$ws.Range("A71").Value = "M3 Lock Nut "
$ws.Range("A72").Value = "M5 Washer "
$ws.Range("B72").Value = "M5 Washer"
$ws.Range("A73").Value = "M3 NUT "
$ws.Range("A74").Value = "M4 Washer "
$ws.Range("B74").Value = "M4 Washer"
$ws.Range("A75").Value = "M3 Washer "
$ws.Range("B75").Value = "M3 Washer"

# --- Hardware and Wheels ---
$ws.Range("A79").Value = "2020 V-Slot 500"
$ws.Range("A80").Value = "2020 V-Slot 480"
$ws.Range("A82").Value = "2028 L Bracket "
$ws.Range("A83").Value = "M5x565 ROD "
$ws.Range("C83").Value = "Belt torque transfer rod"
$ws.Range("A84").Value = "V-Slot Big Wheels"

# --- Printed parts ---
$ws.Range("A117").Value = "M5 Thumb Screw 2 "
$ws.Range("A118").Value = "Y End Stop Buddy "

# Match the author's last-saved cursor position
$ws.Activate()
$ws.Range("B89").Select()
